# Shift the timestamp column (A) forward by 28 days for all data rows (2-97)
# and update the production values (B) for rows 24-42 to reflect the
# retrained model output for the new date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift all timestamps in column A (rows 2 through 97) by +28 days
# (use Value2 -- Value on a date-formatted cell returns a formatted/
# variant representation rather than the underlying numeric serial)
for ($r = 2; $r -le 97; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value2 = $cell.Value2 + 28
}

# Updated "Actual Production (MW)" values for rows 24-42 (column B)
$bValues = @{
    24 = 0
    25 = 0
    26 = 3
    27 = 26
    28 = 76
    29 = 161
    30 = 276
    31 = 399
    32 = 524
    33 = 620
    34 = 789
    35 = 923
    36 = 1011
    37 = 1066
    38 = 1170
    39 = 1266
    40 = 1253
    41 = 1260
    42 = 1284
}

foreach ($r in $bValues.Keys) {
    $ws.Cells.Item($r, 2).Value = $bValues[$r]
}
